# Applet support and SPIF changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "$/pin" header label (F5) must be written before the second
#     table's "c/mm2 =" label so the shared-string table order matches
#     (index 10 = "$/pin", index 11 = "c/mm2 =").
$ws.Range("F5").Value = "$/pin"

# --- New "$/pin" column of formulas on the existing die-size table ---
$ws.Range("F6:F15").Formula = "=E6/A6^2"

# --- Number formats for columns E and F (data rows 6:18) ---
# Style 1: F column -> "0.000"
$ws.Range("F6:F18").NumberFormat = "0.000"
# Style 2: E column -> "0.00"
$ws.Range("E6:E18").NumberFormat = "0.00"

# --- Second (applet / SPIF) table label, row 19 ---
$ws.Range("A19").Value = "c/mm2 ="
# Style 3: right align only (A19 label)
$ws.Range("A19").HorizontalAlignment = -4152

# --- Header row 5 (bold, right aligned), including new F5 label ---
# Style 4: bold + right align (row 5 header cells)
$ws.Range("A5:F5").Font.Bold = $true
$ws.Range("A5:F5").HorizontalAlignment = -4152

# Style 6: left align only (B19 value)
$ws.Range("B19").Value = 14
$ws.Range("B19").HorizontalAlignment = -4131

# --- Column widths (closest values the engine's width quantization allows) ---
$ws.Columns("E").ColumnWidth = 9.6666666666667
$ws.Columns("F").ColumnWidth = 8.1666666666667

# --- Second (applet / SPIF) table body, rows 20:26 ---
$ws.Range("A20").Value = 20
$ws.Range("A21").Value = 32
$ws.Range("A22").Value = 44
$ws.Range("A23").Value = 64
$ws.Range("A24").Value = 88
$ws.Range("A25").Value = 100
$ws.Range("A26").Value = 144

$ws.Range("B20").Formula = "=A20/B$19"
$ws.Range("C20").Formula = "=SQRT(B20)-0.7"
$ws.Range("D20").Formula = "=C20^2"

$ws.Range("B21:B26").Formula = "=A21/B$19"
$ws.Range("C21:C26").Formula = "=SQRT(B21)-0.7"
$ws.Range("D21:D26").Formula = "=C21^2"

# --- Page orientation + final selection/view state ---
$ws.PageSetup.Orientation = 1
$ws.Range("A27").Select()

Write-Host "done"
